$wb = $excel.ActiveWorkbook

# Every "TraceNN" worksheet has an I2 cell that currently holds the text
# "NaN" (a shared string). Excel re-saved the workbook and replaced that
# text value with the numeric literal 0 on every trace sheet.
$traceSheets = @(
    "Trace38","Trace37","Trace36","Trace35","Trace34","Trace33","Trace32",
    "Trace31","Trace30","Trace29","Trace28","Trace27","Trace26","Trace25",
    "Trace24","Trace23","Trace22","Trace21","Trace20","Trace19","Trace18",
    "Trace17","Trace16","Trace15","Trace14","Trace13","Trace12","Trace11",
    "Trace10","Trace9","Trace8","Trace7","Trace6","Trace5","Trace4","Trace3"
)

foreach ($name in $traceSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Cells.Item(2, 9).Value = 0
}

# A handful of sheets also carry refreshed F3/H3 totals (tiny floating
# point drift from the upstream model re-run).
$numericUpdates = @{
    "Trace26" = @{ F3 = 8963712.0121370479; H3 = 9068110.1021370497 }
    "Trace23" = @{ F3 = 8998401.3091335408; H3 = 9103878.9891335424 }
    "Trace14" = @{ F3 = 8526633.0356878955; H3 = 8766551.5236878935 }
    "Trace13" = @{ F3 = 8944580.1798860673; H3 = 9118856.8538860679 }
    "Trace10" = @{ F3 = 8750166.6880306825; H3 = 8938316.5200306848 }
    "Trace6"  = @{ F3 = 8229999.9999991106; H3 = 8515309.8079991098 }
    "Trace32" = @{ F3 = 8855121.830367364;  H3 = 9039071.5503673628 }
}

foreach ($name in $numericUpdates.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = $numericUpdates[$name].F3
    $ws.Range("H3").Value = $numericUpdates[$name].H3
}
